$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (forced as Text to avoid Excel auto-converting
# numeric-looking strings like "1.0000" or "63.00" into numbers, which would lose
# the exact original text representation used in the source workbook.
$updates = @{
    'D2' = '29.454.56'
    'E2' = '  +1.90%  '
    'D3' = '1.854.89'
    'E3' = '  +1.22%  '
    'D4' = '0.9994'
    'E4' = '  -0.06%  '
    'D5' = '244.74'
    'E5' = '  +0.06%  '
    'D6' = '0.6962'
    'E6' = '  +0.89%  '
    'D7' = '1.0000'
    'E7' = '  -0.03%  '
    'D8' = '0.07684'
    'E8' = '  -0.08%  '
    'D9' = '0.3068'
    'E9' = '  +0.65%  '
    'D10' = '23.55'
    'E10' = '  +0.54%  '
    'D11' = '0.07775'
    'E11' = '  -0.50%  '
    'D12' = '5.148'
    'E12' = '  +1.35%  '
    'D13' = '1.857.14'
    'E13' = '  +0.66%  '
    'D14' = '91.01'
    'E14' = '  +0.54%  '
    'D15' = '0.6920'
    'E15' = '  +1.71%  '
    'D16' = '6.263'
    'E16' = '  -2.74%  '
    'D17' = '29.426.72'
    'E17' = '  +1.76%  '
    'D18' = '0.000008341'
    'E18' = '  +0.56%  '
    'D19' = '2.093.09'
    'E19' = '  +0.73%  '
    'D20' = '238.22'
    'E20' = '  -1.94%  '
    'E21' = '  -0.05%  '
    'D22' = '0.9996'
    'E22' = '  -0.01%  '
    'D23' = '7.602'
    'E23' = '  +1.70%  '
    'D24' = '0.9999'
    'E24' = '  -0.03%  '
    'D25' = '0.1493'
    'E25' = '  +1.88%  '
    'D26' = '159.90'
    'E26' = '  -1.89%  '
    'D27' = '8.881'
    'E27' = '  +0.92%  '
    'D28' = '18.24'
    'E28' = '  +0.16%  '
    'D29' = '1.529'
    'E29' = '  -0.96%  '
    'D30' = '4.239'
    'E30' = '  +0.77%  '
    'E31' = '  -0.01%  '
    'D32' = '1.202'
    'E32' = '  +1.64%  '
    'D33' = '0.05091'
    'E33' = '  -0.24%  '
    'D34' = '0.7737'
    'E34' = '  +1.09%  '
    'E35' = '  +2.48%  '
    'D36' = '1.148'
    'E36' = '  +0.54%  '
    'E37' = '  -0.21%  '
    'D38' = '1.317.52'
    'E38' = '  +8.05%  '
    'E39' = '  +1.29%  '
    'D40' = '2.719'
    'E40' = '  +0.80%  '
    'D41' = '0.9475'
    'E41' = '  +0.77%  '
    'D42' = '106.30'
    'E42' = '  -1.57%  '
    'D43' = '5.771'
    'E43' = '  +1.66%  '
    'E44' = '  +0.10%  '
    'D45' = '9.766'
    'E45' = '  +2.24%  '
    'B46' = 'RocketPoolETH'
    'C46' = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
    'D46' = '1.997.81'
    'E46' = '  +1.00%  '
    'B47' = 'Mantle'
    'C47' = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
    'D47' = '0.5217'
    'E47' = '  +0.85%  '
    'B48' = 'RenderToken'
    'C48' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D48' = '1.786'
    'E48' = '  +2.26%  '
    'B49' = 'Aave'
    'C49' = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    'D49' = '63.00'
    'E49' = '  -1.93%  '
    'B50' = 'Aptos'
    'C50' = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    'D50' = '6.952'
    'E50' = '  +0.78%  '
    'B51' = 'Cronos'
    'C51' = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    'D51' = '0.05926'
    'E51' = '  +0.74%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
